# Update the GBIFTaxa worksheet (row 8) so that the example file's GBIF
# taxonomy passes validation: replace the old "Microcopris hidakai" taxon
# row with a new "Morus rubra" (gannet) taxon row, and fill in the
# previously-empty Ignore ID / Parent ID / Comments cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GBIFTaxa")

# Order matters here: new shared strings are interned in the order they
# are first written, so write B8 before A8 to match the expected shared
# string table ordering (Morus rubra, then new_gannet).
$ws.Range("B8").Value = "Morus rubra"
$ws.Range("A8").Value = "new_gannet"

$ws.Range("F8").Value = "Morus"

# The Ignore ID cell (E8) loses its style when its value is replaced.
$ws.Range("E8").ClearFormats()
$ws.Range("E8").Value = 5361886

$ws.Range("H8").Value = 2480962
$ws.Range("I8").Value = "New gannet (not mulberry) species"

# Widen the new Comments column to fit its content and select the cell
# below it, matching the worksheet's final on-screen state.
$ws.Columns.Item(9).ColumnWidth = 29.67
$ws.Range("I14").Select() | Out-Null
